$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values are plain text in the source data (e.g. "26.108.61")
# and some look numeric (e.g. "209.35"). Force text storage via a temporary
# "@" (text) number format, then clear formats again so no style/format diff
# is introduced versus the original (unstyled) data cells.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.108.61"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.72%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.666.76"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.43%  "

$ws.Range("E4").Value = "  -0.39%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "209.35"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -3.92%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5253"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.95%  "

$ws.Range("E7").Value = "  -0.41%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2613"
$ws.Range("D8").ClearFormats()

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06291"
$ws.Range("D9").ClearFormats()

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.09"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.09%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07518"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.23%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.664.32"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.96%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.425"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.12%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5491"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -5.30%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "66.28"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.96%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000007949"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -5.11%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.124.79"
$ws.Range("D17").ClearFormats()

$ws.Range("E18").Value = "  -0.46%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.700"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -4.20%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "186.38"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -4.26%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.24"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -5.72%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.163"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.83%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.003"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.37%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "149.66"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.54%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1241"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -3.71%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.452"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -5.35%  "

$ws.Range("E27").Value = "  +0.03%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06371"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +4.21%  "

$ws.Range("E29").Value = "  -2.76%  "

$ws.Range("E30").Value = "  -4.10%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.489"
$ws.Range("D31").ClearFormats()

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.408"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -4.86%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.634"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.22%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.002"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -3.05%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.404"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.92%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.5997"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -3.23%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.728"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.27%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.106.88"
$ws.Range("D38").ClearFormats()

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.095"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.57%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01613"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.87%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8701"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.07%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.83"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.15%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.817.51"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.36%  "

$ws.Range("E45").Value = "  -2.65%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "55.23"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -4.35%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.003"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.78%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.033"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.54%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05227"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.21%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4245"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.07%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.923"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.20%  "
